$p = $ppt.ActivePresentation

# Slide 22: Title "Bài tập" -> "Bài tập 8.1"
# Split the single run "Bài tập" after "Bài " (chars 5-7 = "tập") and
# append " 8.1" to the second run's text.
$s22 = $p.Slides.Item(22)
$tr22 = $s22.Shapes.Item(1).TextFrame.TextRange
$sub22 = $tr22.Characters(5, 3)
$sub22.Text = "tập 8.1"

# Slide 23: Title "Bài tập 2" -> "Bài tập 8.2"
# Replace the trailing "2" (last character) with "8.2"; this splits the
# run containing " 2" into a " " run and an "8.2" run.
$s23 = $p.Slides.Item(23)
$tr23 = $s23.Shapes.Item(1).TextFrame.TextRange
$sub23 = $tr23.Characters($tr23.Length, 1)
$sub23.Text = "8.2"
